# Applies the cryptos list update (prices, % changes, and a 3-row reorder
# for THORChain / Stellar / ApeXProtocol) captured by the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows 49-51: Coin name & Link cells (text, safe to assign directly) ---
$ws.Range('B49').Value = 'THORChain'
$ws.Range('C49').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('B50').Value = 'Stellar'
$ws.Range('C50').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('B51').Value = 'ApeXProtocol'
$ws.Range('C51').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'

# --- Price column (D): force text format first so values such as "1.00",
#     "0.430", "8.38" etc. keep their exact textual representation instead
#     of being auto-converted to numbers by Excel. ---
$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '67.819.23'
$ws.Range('D2').Style = "Normal"
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '3.601.47'
$ws.Range('D3').Style = "Normal"
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '1.00'
$ws.Range('D4').Style = "Normal"
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '202.63'
$ws.Range('D5').Style = "Normal"
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '563.48'
$ws.Range('D6').Style = "Normal"
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '3.599.96'
$ws.Range('D7').Style = "Normal"
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.672'
$ws.Range('D10').Style = "Normal"
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '60.55'
$ws.Range('D11').Style = "Normal"
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.151'
$ws.Range('D12').Style = "Normal"
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '0.0000287'
$ws.Range('D13').Style = "Normal"
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '10.01'
$ws.Range('D14').Style = "Normal"
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '4.201.36'
$ws.Range('D15').Style = "Normal"
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '3.614.44'
$ws.Range('D16').Style = "Normal"
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '18.92'
$ws.Range('D18').Style = "Normal"
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '67.823.16'
$ws.Range('D19').Style = "Normal"
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '12.33'
$ws.Range('D20').Style = "Normal"
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '403.26'
$ws.Range('D22').Style = "Normal"
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '12.97'
$ws.Range('D23').Style = "Normal"
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '85.17'
$ws.Range('D25').Style = "Normal"
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '2.95'
$ws.Range('D26').Style = "Normal"
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '6.11'
$ws.Range('D29').Style = "Normal"
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '8.38'
$ws.Range('D30').Style = "Normal"
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '31.54'
$ws.Range('D32').Style = "Normal"
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '682.95'
$ws.Range('D33').Style = "Normal"
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '12.19'
$ws.Range('D34').Style = "Normal"
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '42.33'
$ws.Range('D37').Style = "Normal"
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.430'
$ws.Range('D38').Style = "Normal"
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.999'
$ws.Range('D39').Style = "Normal"
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.0₃0773'
$ws.Range('D40').Style = "Normal"
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '3.235.63'
$ws.Range('D42').Style = "Normal"
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '2.79'
$ws.Range('D44').Style = "Normal"
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '3.01'
$ws.Range('D45').Style = "Normal"
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '0.999'
$ws.Range('D46').Style = "Normal"
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '0.0418'
$ws.Range('D47').Style = "Normal"
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '8.83'
$ws.Range('D49').Style = "Normal"
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '0.131'
$ws.Range('D50').Style = "Normal"
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '3.10'
$ws.Range('D51').Style = "Normal"

# --- Volume(1h) column (E): percentage text with surrounding spaces ---
$ws.Range('E2').Value = '  +2.14%  '
$ws.Range('E3').Value = '  +0.82%  '
$ws.Range('E4').Value = '  +0.06%  '
$ws.Range('E5').Value = '  +8.12%  '
$ws.Range('E6').Value = '  -4.27%  '
$ws.Range('E7').Value = '  +0.94%  '
$ws.Range('E8').Value = '  +0.26%  '
$ws.Range('E10').Value = '  +0.20%  '
$ws.Range('E11').Value = '  +12.87%  '
$ws.Range('E12').Value = '  +3.22%  '
$ws.Range('E13').Value = '  +10.24%  '
$ws.Range('E14').Value = '  +2.27%  '
$ws.Range('E15').Value = '  +1.26%  '
$ws.Range('E16').Value = '  +1.20%  '
$ws.Range('E17').Value = '  +0.53%  '
$ws.Range('E18').Value = '  +3.35%  '
$ws.Range('E19').Value = '  +2.07%  '
$ws.Range('E20').Value = '  +0.55%  '
$ws.Range('E21').Value = '  +1.84%  '
$ws.Range('E22').Value = '  +1.70%  '
$ws.Range('E23').Value = '  +13.83%  '
$ws.Range('E24').Value = '  -4.89%  '
$ws.Range('E25').Value = '  -0.74%  '
$ws.Range('E26').Value = '  +2.12%  '
$ws.Range('E27').Value = '  +0.97%  '
$ws.Range('E28').Value = '  +9.26%  '
$ws.Range('E29').Value = '  +1.23%  '
$ws.Range('E30').Value = '  +18.20%  '
$ws.Range('E31').Value = '  +5.20%  '
$ws.Range('E32').Value = '  +1.34%  '
$ws.Range('E33').Value = '  +10.83%  '
$ws.Range('E34').Value = '  +0.17%  '
$ws.Range('E35').Value = '  +0.44%  '
$ws.Range('E36').Value = '  +0.57%  '
$ws.Range('E37').Value = '  +2.38%  '
$ws.Range('E38').Value = '  +10.04%  '
$ws.Range('E39').Value = '  -0.07%  '
$ws.Range('E40').Value = '  +1.98%  '
$ws.Range('E41').Value = '  +15.12%  '
$ws.Range('E42').Value = '  +6.46%  '
$ws.Range('E43').Value = '  +3.32%  '
$ws.Range('E44').Value = '  +10.18%  '
$ws.Range('E45').Value = '  +28.82%  '
$ws.Range('E46').Value = '  -0.03%  '
$ws.Range('E47').Value = '  +2.14%  '
$ws.Range('E48').Value = '  +9.87%  '
$ws.Range('E49').Value = '  +2.71%  '
$ws.Range('E50').Value = '  +0.33%  '
$ws.Range('E51').Value = '  +3.48%  '
